$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$lines = @(
    "Physical level:",
    "TS data require different data layout than graph data;",
    "LSM-Tree-like (e.g., RocksDB)",
    "InfluxDB 3.0 on Parquet.",
    "Metadata modelling (ausiliary structures);",
    "query formalization and optimization;",
    "Analytics:",
    "TS operators in Cypher/GQL (Graph analytics);",
    "shape/patthern matching;",
    "Cross time-series operators:",
    "Identify plants/grids with similar drying patterns over the last 24h",
    "Graph-TS cross-operators:",
    "Correlate graph metrics (node degree, node/edge properties) with time-series trends",
    "e.g., landslide monitoring sensor network: correlation between pressure measurements and dynamic edge weights between nearby sensors",
    "Correlate soil drying with temperature (spatial join with ARPAE weather stations)",
    "LLMs:",
    "Text to query (in hybrid models)",
    "Repair"
)

$levels = @(0, 1, 2, 2, 1, 1, 0, 1, 2, 1, 2, 1, 2, 3, 2, 0, 1, 1)

# Setting TextRange.Text directly to a string with the same paragraph count
# as the existing text makes the host diff new vs. old paragraph-by-paragraph
# (preserving common prefixes/suffixes as separate runs). Break that paragraph
# -count match first with an unrelated placeholder, then apply the real text
# so every paragraph ends up as a single clean run.
$tr.Text = "placeholder"
$tr.Text = [string]::Join("`r", $lines)

for ($i = 1; $i -le $lines.Count; $i++) {
    $sub = $tr.Paragraphs($i, 1)
    $sub.IndentLevel = $levels[$i - 1] + 1
}
